# Apply the commit's changes to the workbook:
#  - Update the "Date" metadata value on the Metadata sheet
#  - Change the Relationship value for the "bmi" source mapping rows
#    (in the LOINC and UMLS mapping tables) from "equivalent" to
#    "source-is-narrower-than-target"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value (row 8, column B) ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2022-04-11T21:19:38+00:00"

# --- Mapping Table 0 (UMLS): bmi -> C1305855 row, Relationship column C, row 6 ---
$mapping0 = $wb.Worksheets.Item("Mapping Table 0")
$mapping0.Range("C6").Value = "source-is-narrower-than-target"

# --- Mapping Table 1 (LOINC): bmi -> LP35925-4 row, Relationship column C, row 3 ---
$mapping1 = $wb.Worksheets.Item("Mapping Table 1")
$mapping1.Range("C3").Value = "source-is-narrower-than-target"
